$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 (Sheet1 -> TC01_LoginTest): add a new column C mirroring "username"
# header, plus new credential values in C2/C3. Order of writes matters for
# shared-string allocation (pradeep before test).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("C1").Value = "username"
$ws1.Range("B3").Value = "pradeep5"
$ws1.Range("C3").Value = "pradeep"
$ws1.Range("B2").Value = "pradeep5"
$ws1.Range("C2").Value = "test"

# Cell A2 picks up the Hyperlink cell style (s="1") in the new layout.
$ws1.Range("A2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Sheet2 (Sheet2 -> TC02_SignUp): replace sample sign-up values, widen the
# columns, drop the mailto hyperlink (keep the Hyperlink cell style on B2).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Hyperlinks.Delete()
$ws2.Range("C2").Value = "testuser123"
$ws2.Range("B2").Value = "Automationuser"
$ws2.Range("A2").Value = "AutomationUser"

$ws2.Columns.Item(1).ColumnWidth = 17.333333333333332
$ws2.Columns.Item(2).ColumnWidth = 16.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# Sheet3 (Sheet3 -> TC03_CreateTeam): fill in the retro-template name column.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A1").Value = "template"
$ws3.Range("A2").Value = " Liked, Learned, Lacked, Longed for "
$ws3.Range("A3").Value = " Anchors, Engines "
$ws3.Range("A4").Value = " Drop,Add,Keep,Improve "
$ws3.Range("A5").Value = " Future considerations,Lessons learned, Accomplishments,Problem areas "
$ws3.Range("A6").Value = " Liked, Learned, Lacked, Longed for "
$ws3.Range("A7").Value = " Mad, Sad, Glad "
$ws3.Range("A8").Value = " Start, More of, Continue, Less of, Stop "
$ws3.Range("A9").Value = " Start, Stop "
$ws3.Range("A10").Value = " Start, Stop, Continue "
$ws3.Range("A11").Value = " Wishes, Risks, Appreciations, Puzzles "
$ws3.Range("A12").Value = " Design your own "

$ws3.Columns.Item(1).ColumnWidth = 31.666666666666668

# ---------------------------------------------------------------------------
# Sheet names.
# ---------------------------------------------------------------------------
$ws1.Name = "TC01_LoginTest"
$ws2.Name = "TC02_SignUp"
$ws3.Name = "TC03_CreateTeam"

# ---------------------------------------------------------------------------
# Selections / active sheet: TC03_CreateTeam becomes the active tab.
# ---------------------------------------------------------------------------
$ws1.Range("C9").Select() | Out-Null
$ws2.Range("A3").Select() | Out-Null
$ws3.Select() | Out-Null
$ws3.Range("I16").Select() | Out-Null

# Page orientation on sheet1.
$ws1.PageSetup.Orientation = 1
